# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates to Sheets/Atomos_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: ALC
$ws.Cells.Item(6, 8).Value = 1329.3529
$ws.Cells.Item(6, 9).Value = 204
$ws.Cells.Item(6, 11).Value = 612
$ws.Cells.Item(6, 13).Value = -500

# Row 9: ALC
$ws.Cells.Item(9, 8).Value = 7692592
$ws.Cells.Item(9, 10).Value = 464.2857
$ws.Cells.Item(9, 12).Value = 464.2857
$ws.Cells.Item(9, 14).Value = -802.2857

# Row 28: ALC
$ws.Cells.Item(28, 8).Value = 439.94116
$ws.Cells.Item(28, 9).Value = 430.2143
$ws.Cells.Item(28, 11).Value = 430.2143
$ws.Cells.Item(28, 13).Value = 54.78570000000002

# Row 38: ALC
$ws.Cells.Item(38, 8).Value = 357.7143
$ws.Cells.Item(38, 9).Value = 259.85184
$ws.Cells.Item(38, 11).Value = 779.5555199999999
$ws.Cells.Item(38, 13).Value = -407.5555199999999

# Row 107: ALC
$ws.Cells.Item(107, 8).Value = 519.26666
$ws.Cells.Item(107, 9).Value = 407.41666
$ws.Cells.Item(107, 10).Value = 966.6667
$ws.Cells.Item(107, 11).Value = 407.41666
$ws.Cells.Item(107, 12).Value = 966.6667
$ws.Cells.Item(107, 13).Value = 1512.58334
$ws.Cells.Item(107, 14).Value = -4806.6667

# Row 112: ALC
$ws.Cells.Item(112, 8).Value = 1733.5714
$ws.Cells.Item(112, 10).Value = 1240
$ws.Cells.Item(112, 12).Value = 3720
$ws.Cells.Item(112, 14).Value = -5936

# Row 134: ALC
$ws.Cells.Item(134, 8).Value = 27595.715
$ws.Cells.Item(134, 10).Value = 27595.715
$ws.Cells.Item(134, 12).Value = 27595.715
$ws.Cells.Item(134, 14).Value = -37735.715

# Row 136: ALC
$ws.Cells.Item(136, 8).Value = 26736.334
$ws.Cells.Item(136, 9).Value = 20709
$ws.Cells.Item(136, 10).Value = 29750
$ws.Cells.Item(136, 11).Value = 20709
$ws.Cells.Item(136, 12).Value = 29750
$ws.Cells.Item(136, 13).Value = -15609
$ws.Cells.Item(136, 14).Value = -39950

# Row 139: ALC
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()  # N139 removed (was -39280)

# Row 141: ALC
$ws.Cells.Item(141, 8).Value = 688052.9
$ws.Cells.Item(141, 9).Value = 1670.2727
$ws.Cells.Item(141, 10).Value = 1766654.1
$ws.Cells.Item(141, 11).Value = 5010.8181
$ws.Cells.Item(141, 12).Value = 5299962.300000001
$ws.Cells.Item(141, 13).Value = 169.1818999999996
$ws.Cells.Item(141, 14).Value = -5310322.300000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: ARM
$ws.Cells.Item(32, 8).Value = 17757.6
$ws.Cells.Item(32, 9).Value = 14373.24
$ws.Cells.Item(32, 10).Value = 29038.8
$ws.Cells.Item(32, 11).Value = 14373.24
$ws.Cells.Item(32, 12).Value = 29038.8
$ws.Cells.Item(32, 13).Value = -14086.24
$ws.Cells.Item(32, 14).Value = -29612.8

# Row 132: ARM
$ws.Cells.Item(132, 8).Value = 2504.8918
$ws.Cells.Item(132, 9).Value = 1865.76
$ws.Cells.Item(132, 10).Value = 3836.4167
$ws.Cells.Item(132, 11).Value = 5597.28
$ws.Cells.Item(132, 12).Value = 11509.2501
$ws.Cells.Item(132, 13).Value = -3067.28
$ws.Cells.Item(132, 14).Value = -16569.2501

# Row 135: ARM
$ws.Cells.Item(135, 8).Value = 39800
$ws.Cells.Item(135, 10).Value = 39800
$ws.Cells.Item(135, 12).Value = 39800
$ws.Cells.Item(135, 14).Value = -49940

# Row 137: ARM
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()  # N137 removed (was -40200)

# Row 139: ARM
$ws.Cells.Item(139, 8).Value = 39800
$ws.Cells.Item(139, 10).Value = 39800
$ws.Cells.Item(139, 12).Value = 39800
$ws.Cells.Item(139, 14).Value = -50080

$ws = $wb.Worksheets.Item("CRP")
# Row 141: CRP
$ws.Cells.Item(141, 8).Value = 29965.518
$ws.Cells.Item(141, 9).Value = 4333.3335
$ws.Cells.Item(141, 10).Value = 32923.08
$ws.Cells.Item(141, 11).Value = 4333.3335
$ws.Cells.Item(141, 12).Value = 32923.08
$ws.Cells.Item(141, 13).Value = 846.6665000000003
$ws.Cells.Item(141, 14).Value = -43283.08

$ws = $wb.Worksheets.Item("CUL")
# Row 131: CUL
$ws.Cells.Item(131, 8).Value = 1522.6041
$ws.Cells.Item(131, 9).Value = 2364.0908
$ws.Cells.Item(131, 10).Value = 1272.4324
$ws.Cells.Item(131, 11).Value = 7092.2724
$ws.Cells.Item(131, 12).Value = 3817.2972
$ws.Cells.Item(131, 13).Value = -2052.2724
$ws.Cells.Item(131, 14).Value = -13897.2972

# Row 133: CUL
$ws.Cells.Item(133, 8).Value = 5420.8823
$ws.Cells.Item(133, 9).Value = 4998.3335
$ws.Cells.Item(133, 10).Value = 5651.364
$ws.Cells.Item(133, 11).Value = 14995.0005
$ws.Cells.Item(133, 12).Value = 16954.092
$ws.Cells.Item(133, 13).Value = -9935.000499999998
$ws.Cells.Item(133, 14).Value = -27074.092

# Row 137: CUL
$ws.Cells.Item(137, 8).Value = 3104.4814
$ws.Cells.Item(137, 9).Value = 2493.889
$ws.Cells.Item(137, 10).Value = 4325.6665
$ws.Cells.Item(137, 11).Value = 7481.667
$ws.Cells.Item(137, 12).Value = 12976.9995
$ws.Cells.Item(137, 13).Value = -2381.667
$ws.Cells.Item(137, 14).Value = -23176.9995

# Row 139: CUL
$ws.Cells.Item(139, 8).Value = 11367686
$ws.Cells.Item(139, 9).Value = 13160957
$ws.Cells.Item(139, 10).Value = 10306.667
$ws.Cells.Item(139, 11).Value = 39482871
$ws.Cells.Item(139, 12).Value = 30920.001
$ws.Cells.Item(139, 13).Value = -39477731
$ws.Cells.Item(139, 14).Value = -41200.001

$ws = $wb.Worksheets.Item("GSM")
# Row 11: GSM
$ws.Cells.Item(11, 8).Value = 4803000
$ws.Cells.Item(11, 9).Value = 7114444.5
$ws.Cells.Item(11, 10).Value = 2911818.2
$ws.Cells.Item(11, 11).Value = 7114444.5
$ws.Cells.Item(11, 12).Value = 2911818.2
$ws.Cells.Item(11, 13).Value = -7114305.5
$ws.Cells.Item(11, 14).Value = -2912096.2

# Row 18: GSM
$ws.Cells.Item(18, 8).Value = 53576.57
$ws.Cells.Item(18, 10).Value = 53576.57
$ws.Cells.Item(18, 12).Value = 53576.57
$ws.Cells.Item(18, 14).Value = -54162.57

# Row 57: GSM
$ws.Cells.Item(57, 8).Value = 19630.5
$ws.Cells.Item(57, 10).Value = 23695.75
$ws.Cells.Item(57, 12).Value = 23695.75
$ws.Cells.Item(57, 14).Value = -25335.75

# Row 88: GSM
$ws.Cells.Item(88, 8).Value = 29333.334
$ws.Cells.Item(88, 10).Value = 29333.334
$ws.Cells.Item(88, 12).Value = 29333.334
$ws.Cells.Item(88, 14).Value = -30235.334

# Row 91: GSM
$ws.Cells.Item(91, 8).Value = 29333.334
$ws.Cells.Item(91, 10).Value = 29333.334
$ws.Cells.Item(91, 12).Value = 29333.334
$ws.Cells.Item(91, 14).Value = -32453.334

$ws = $wb.Worksheets.Item("LTW")
# Row 20: LTW
$ws.Cells.Item(20, 8).Value = 44291.715
$ws.Cells.Item(20, 10).Value = 44291.715
$ws.Cells.Item(20, 12).Value = 44291.715
$ws.Cells.Item(20, 14).Value = -44743.715

# Row 25: LTW
$ws.Cells.Item(25, 8).Value = 40008
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 40008
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 40008
$ws.Cells.Item(25, 13).ClearContents()  # M25 removed (was -24769.5)
$ws.Cells.Item(25, 14).Value = -40468

# Row 61: LTW
$ws.Cells.Item(61, 8).Value = 166672770
$ws.Cells.Item(61, 9).Value = 500000500
$ws.Cells.Item(61, 10).Value = 8900
$ws.Cells.Item(61, 11).Value = 500000500
$ws.Cells.Item(61, 12).Value = 8900
$ws.Cells.Item(61, 13).Value = -500000298
$ws.Cells.Item(61, 14).Value = -9304

# Row 93: LTW
$ws.Cells.Item(93, 8).Value = 2910.889
$ws.Cells.Item(93, 9).Value = 2149.8333
$ws.Cells.Item(93, 10).Value = 4433
$ws.Cells.Item(93, 11).Value = 2149.8333
$ws.Cells.Item(93, 12).Value = 4433
$ws.Cells.Item(93, 13).Value = -901.8332999999998
$ws.Cells.Item(93, 14).Value = -6929

# Row 113: LTW
$ws.Cells.Item(113, 8).Value = 166672770
$ws.Cells.Item(113, 9).Value = 500000500
$ws.Cells.Item(113, 10).Value = 8900
$ws.Cells.Item(113, 11).Value = 500000500
$ws.Cells.Item(113, 12).Value = 8900
$ws.Cells.Item(113, 13).Value = -499998330
$ws.Cells.Item(113, 14).Value = -13240

$ws = $wb.Worksheets.Item("WVR")
# Row 49: WVR
$ws.Cells.Item(49, 8).Value = 3750
$ws.Cells.Item(49, 10).Value = 3750
$ws.Cells.Item(49, 12).Value = 3750
$ws.Cells.Item(49, 14).Value = -4210

# Row 54: WVR
$ws.Cells.Item(54, 8).Value = 21751.334
$ws.Cells.Item(54, 10).Value = 21751.334
$ws.Cells.Item(54, 12).Value = 21751.334
$ws.Cells.Item(54, 14).Value = -22791.334

# Row 113: WVR
$ws.Cells.Item(113, 8).Value = 5140
$ws.Cells.Item(113, 9).Value = 280
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 11).Value = 840
$ws.Cells.Item(113, 12).Value = 30000
$ws.Cells.Item(113, 13).Value = 1330
$ws.Cells.Item(113, 14).Value = -34340

# Row 119: WVR
$ws.Cells.Item(119, 8).Value = 29666.666
$ws.Cells.Item(119, 10).Value = 29666.666
$ws.Cells.Item(119, 12).Value = 29666.666
$ws.Cells.Item(119, 14).Value = -39342.666

# Row 135: WVR
$ws.Cells.Item(135, 8).Value = 75102.5
$ws.Cells.Item(135, 10).Value = 75102.5
$ws.Cells.Item(135, 12).Value = 75102.5
$ws.Cells.Item(135, 14).Value = -85242.5

# Row 137: WVR
$ws.Cells.Item(137, 8).Value = 54266.668
$ws.Cells.Item(137, 10).Value = 66900
$ws.Cells.Item(137, 12).Value = 66900
$ws.Cells.Item(137, 14).Value = -77100

# Row 139: WVR
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()  # N139 removed (was -40280)

# Row 141: WVR
$ws.Cells.Item(141, 8).Value = 28850
$ws.Cells.Item(141, 10).Value = 28850
$ws.Cells.Item(141, 12).Value = 28850
$ws.Cells.Item(141, 14).Value = -39210
